$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 793, shifting existing rows 793..821 down to 794..822.
# This preserves all the existing data/formatting for rows that move down,
# matching the diff where each old row N's content becomes row N+1's content.
$ws.Rows.Item(793).Insert()

# Populate the newly inserted row 793 with its new values (per diff).
$ws.Range("A793").Value2 = 5
$ws.Range("B793").Value2 = "Macroferia Regional de Talca"
$ws.Range("C793").Value2 = "Maule"
$ws.Range("D793").Value2 = 45075
$ws.Range("D793").NumberFormat = $ws.Range("D794").NumberFormat
$ws.Range("E793").Value2 = 7
$ws.Range("F793").Value2 = "Fruta"
$ws.Range("G793").Value2 = 100101
$ws.Range("H793").Value2 = "Berries"
$ws.Range("I793").Value2 = 100112025
$ws.Range("J793").Value2 = "Frutilla"
$ws.Range("K793").Value2 = "Sin especificar"
$ws.Range("L793").Value2 = "Primera"
$ws.Range("M793").Value2 = 180
$ws.Range("N793").Value2 = 13000
$ws.Range("O793").Value2 = 13000
$ws.Range("P793").Value2 = 13000
$ws.Range("Q793").Value2 = "$/bandeja 7 kilos"
$ws.Range("R793").Value2 = "Provincia de Melipilla"
$ws.Range("S793").Value2 = 1857
$ws.Range("T793").Value2 = 7
